$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-07-03 Wednesday"

# The worksheet body is a single 20-row x 5-column table; every 4th row
# (1, 5, 9, 13, 17) holds the division problems, the rows between are
# blank answer rows. Address each problem cell positionally (row, col)
# rather than by text match, since several old/new values repeat
# elsewhere in the table and a left-to-right Find/Replace would risk
# matching a cell that was already updated.
$t = $d.Tables.Item(1)

$updates = @(
    @{r=1;  c=1; v="71÷6="},
    @{r=1;  c=2; v="66÷2="},
    @{r=1;  c=3; v="52÷4="},
    @{r=1;  c=4; v="16÷2="},
    @{r=1;  c=5; v="52÷3="},

    @{r=5;  c=1; v="63÷9="},
    @{r=5;  c=2; v="20÷9="},
    @{r=5;  c=3; v="95÷6="},
    @{r=5;  c=4; v="35÷2="},
    @{r=5;  c=5; v="30÷9="},

    @{r=9;  c=1; v="68÷4="},
    @{r=9;  c=2; v="65÷7="},
    @{r=9;  c=3; v="69÷3="},
    @{r=9;  c=4; v="38÷2="},
    @{r=9;  c=5; v="43÷7="},

    @{r=13; c=1; v="51÷2="},
    @{r=13; c=2; v="89÷7="},
    @{r=13; c=3; v="86÷2="},
    @{r=13; c=4; v="88÷7="},
    @{r=13; c=5; v="38÷7="},

    @{r=17; c=1; v="78÷8="},
    @{r=17; c=2; v="42÷9="},
    @{r=17; c=3; v="68÷4="},
    @{r=17; c=4; v="23÷6="},
    @{r=17; c=5; v="95÷7="}
)

foreach ($u in $updates) {
    $t.Cell($u.r, $u.c).Range.Text = $u.v
}

Write-Host "applied" $updates.Count "cell updates"
